# Update Work Week and Social Spending
# Updates the "GDP per Capita" data series on the "Data" sheet for Comoros:
#   - Rows 2-60  (years 1950-2008): replace the existing value with the
#     revised figure from the updated source dataset.
#   - Rows 61-68 (years 2009-2016): brand new rows appended to extend the
#     series through 2016.
# All values in column E are text (shared-string) cells, matching the
# original file's storage, so we force a "Text" number format before
# writing each value and then clear the formatting again so no stray
# cell styles are left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New GDP per Capita values for years 1950 (row 2) .. 2016 (row 68), in order.
$gdpValues = @(
    "534",
    "553",
    "559",
    "571",
    "590",
    "596",
    "606",
    "615",
    "618",
    "639",
    "679",
    "671",
    "714",
    "846",
    "888",
    "870",
    "937",
    "952",
    "928",
    "918",
    "961",
    "974",
    "966",
    "1006",
    "1004",
    "1038",
    "1009",
    "966",
    "964",
    "1023",
    "1065",
    "1081",
    "1124",
    "1149",
    "1168",
    "1162",
    "1157",
    "1149",
    "1152",
    "1108",
    "1092",
    "1101.00358737775",
    "1185.51412091475",
    "1210.63529556643",
    "1136.74353640223",
    "1166.93684731292",
    "1141.60158685115",
    "1177.58252420305",
    "1179.84405666537",
    "1189.2566885002",
    "1204.40239235915",
    "1217.49922137723",
    "1229.85592357994",
    "1240.05133670709",
    "1247.27494837471",
    "1265.68670368259",
    "1388.00263430292",
    "1382.53908934196",
    "1419.36709608481",
    "1455.28233833545",
    "1432.56242087922",
    "1479",
    "1509",
    "1619",
    "1647",
    "1696",
    "1702"
)

$firstRow = 2
$lastRow = $firstRow + $gdpValues.Length - 1   # 68

# Make sure the whole destination range is stored as text so the
# numeric-looking strings aren't silently converted to numbers.
$eRange = $ws.Range("E$firstRow`:E$lastRow")
$eRange.NumberFormat = "@"

for ($i = 0; $i -lt $gdpValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $gdpValues[$i]
}

# Drop the temporary "Text" number format so no extra cell style sticks
# around on cells that didn't have one before.
$eRange.ClearFormats()

# Fill in columns A-D for the 8 newly appended rows (2009-2016).
$newRowsFirst = 61
$newRowsLast = 68
$year = 2009
for ($row = $newRowsFirst; $row -le $newRowsLast; $row++) {
    $ws.Cells.Item($row, 1).Value = 174
    $ws.Cells.Item($row, 2).Value = "Comoros"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $year = $year + 1
}

Write-Output "Updated E$firstRow`:E$lastRow and appended rows $newRowsFirst-$newRowsLast"
